# Se procesan de nuevo los datos con las nuevas dimensiones curadas
# Update the "municipio-nombre" (column E) and "contenedores-de-pilas"
# (column I) metadata columns to reflect the newly curated dimensions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E: municipio-nombre
$ws.Range("E2").Value = "sdmx-dimension:refArea"
$ws.Range("E3").Value = "dim"
$ws.Range("E4").Value = "URI-Municipio"

# Column I: contenedores-de-pilas
$ws.Range("I2").Value = "iaest-measure:contenedores-de-pilas"
$ws.Range("I3").Value = "medida"
$ws.Range("I4").Value = "xsd:int"
$ws.Range("I5").Clear()
